$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value2 = 2
$ws.Range("G2").Value2 = 40.7349555
$ws.Range("H2").Value2 = 81.469911
$ws.Range("I2").Value2 = 0.05567871843833241
$ws.Range("J2").Value2 = 0.03826666865920979
$ws.Range("K2").Value2 = 2
$ws.Range("M2").Value2 = 14.989415
$ws.Range("N2").Value2 = 29.97883
$ws.Range("O2").Value2 = 0.05547446260572933
$ws.Range("P2").Value2 = 0.03893791130463959
$ws.Range("Q2").Value2 = 610.5931529960325
$ws.Range("R2").Value2 = 2442.37261198413
$ws.Range("S2").Value2 = 0.003088746983942204
$ws.Range("T2").Value2 = 0.001490024150176342
$ws.Range("E3").Value2 = 2
$ws.Range("G3").Value2 = 40.7349555
$ws.Range("H3").Value2 = 81.469911
$ws.Range("I3").Value2 = 0.05567871843833241
$ws.Range("J3").Value2 = 0.03826666865920979
$ws.Range("K3").Value2 = 3
$ws.Range("M3").Value2 = 58.255493
$ws.Range("N3").Value2 = 174.766479
$ws.Range("O3").Value2 = 0.2155982850569436
$ws.Range("P3").Value2 = 0.2269949046819425
$ws.Range("Q3").Value2 = 2373.034914985561
$ws.Range("R3").Value2 = 14238.20948991337
$ws.Range("S3").Value2 = 0.01200423620947289
$ws.Range("T3").Value2 = 0.008686338804792802
$ws.Range("E4").Value2 = 2
$ws.Range("G4").Value2 = 40.7349555
$ws.Range("H4").Value2 = 81.469911
$ws.Range("I4").Value2 = 0.05567871843833241
$ws.Range("J4").Value2 = 0.03826666865920979
$ws.Range("K4").Value2 = 3
$ws.Range("M4").Value2 = 80.178917
$ws.Range("N4").Value2 = 240.536751
$ws.Range("O4").Value2 = 0.2967348847759819
$ws.Range("P4").Value2 = 0.3124204205415681
$ws.Range("Q4").Value2 = 3266.084616033193
$ws.Range("R4").Value2 = 19596.50769619916
$ws.Range("S4").Value2 = 0.01652181810027291
$ws.Range("T4").Value2 = 0.01195528871523517
$ws.Range("E5").Value2 = 2
$ws.Range("G5").Value2 = 40.7349555
$ws.Range("H5").Value2 = 81.469911
$ws.Range("I5").Value2 = 0.05567871843833241
$ws.Range("J5").Value2 = 0.03826666865920979
$ws.Range("K5").Value2 = 3
$ws.Range("M5").Value2 = 79.34548433333335
$ws.Range("N5").Value2 = 238.036453
$ws.Range("O5").Value2 = 0.2936504262229702
$ws.Range("P5").Value2 = 0.3091729161606711
$ws.Range("Q5").Value2 = 3232.134773444281
$ws.Range("R5").Value2 = 19392.80864066568
$ws.Range("S5").Value2 = 0.01635007940096506
$ws.Range("T5").Value2 = 0.01183101754112205
$ws.Range("E6").Value2 = 2
$ws.Range("G6").Value2 = 40.7349555
$ws.Range("H6").Value2 = 81.469911
$ws.Range("I6").Value2 = 0.05567871843833241
$ws.Range("J6").Value2 = 0.03826666865920979
$ws.Range("K6").Value2 = 3
$ws.Range("M6").Value2 = 11.726012
$ws.Range("N6").Value2 = 35.17803600000001
$ws.Range("O6").Value2 = 0.04339690469630293
$ws.Range("P6").Value2 = 0.04569088405516222
$ws.Range("Q6").Value2 = 477.6585770124661
$ws.Range("R6").Value2 = 2865.951462074796
$ws.Range("S6").Value2 = 0.002416284037680596
$ws.Range("T6").Value2 = 0.001748437920885264
$ws.Range("E7").Value2 = 2
$ws.Range("G7").Value2 = 40.7349555
$ws.Range("H7").Value2 = 81.469911
$ws.Range("I7").Value2 = 0.05567871843833241
$ws.Range("J7").Value2 = 0.03826666865920979
$ws.Range("K7").Value2 = 2
$ws.Range("M7").Value2 = 25.708558
$ws.Range("N7").Value2 = 51.417116
$ws.Range("O7").Value2 = 0.09514503664207198
$ws.Range("P7").Value2 = 0.0667829632560165
$ws.Range("Q7").Value2 = 1047.236966099169
$ws.Range("R7").Value2 = 4188.947864396676
$ws.Range("S7").Value2 = 0.005297553705998746
$ws.Range("T7").Value2 = 0.002555561526998166
$ws.Range("E8").Value2 = 3
$ws.Range("G8").Value2 = 350.3919066666667
$ws.Range("H8").Value2 = 1051.17572
$ws.Range("I8").Value2 = 0.4789344206933965
$ws.Range("J8").Value2 = 0.4937404802104949
$ws.Range("K8").Value2 = 2
$ws.Range("M8").Value2 = 14.989415
$ws.Range("N8").Value2 = 29.97883
$ws.Range("O8").Value2 = 0.05547446260572933
$ws.Range("P8").Value2 = 0.03893791130463959
$ws.Range("Q8").Value2 = 5252.169701667934
$ws.Range("R8").Value2 = 31513.0182100076
$ws.Range("S8").Value2 = 0.02656862961135246
$ws.Range("T8").Value2 = 0.01922522302594641
$ws.Range("E9").Value2 = 3
$ws.Range("G9").Value2 = 350.3919066666667
$ws.Range("H9").Value2 = 1051.17572
$ws.Range("I9").Value2 = 0.4789344206933965
$ws.Range("J9").Value2 = 0.4937404802104949
$ws.Range("K9").Value2 = 3
$ws.Range("M9").Value2 = 58.255493
$ws.Range("N9").Value2 = 174.766479
$ws.Range("O9").Value2 = 0.2155982850569436
$ws.Range("P9").Value2 = 0.2269949046819425
$ws.Range("Q9").Value2 = 20412.25326607666
$ws.Range("R9").Value2 = 183710.2793946899
$ws.Range("S9").Value2 = 0.1032574397562371
$ws.Range("T9").Value2 = 0.1120765732429978
$ws.Range("E10").Value2 = 3
$ws.Range("G10").Value2 = 350.3919066666667
$ws.Range("H10").Value2 = 1051.17572
$ws.Range("I10").Value2 = 0.4789344206933965
$ws.Range("J10").Value2 = 0.4937404802104949
$ws.Range("K10").Value2 = 3
$ws.Range("M10").Value2 = 80.178917
$ws.Range("N10").Value2 = 240.536751
$ws.Range("O10").Value2 = 0.2967348847759819
$ws.Range("P10").Value2 = 0.3124204205415681
$ws.Range("Q10").Value2 = 28094.04360209841
$ws.Range("R10").Value2 = 252846.3924188857
$ws.Range("S10").Value2 = 0.1421165501397066
$ws.Range("T10").Value2 = 0.1542546084657586
$ws.Range("E11").Value2 = 3
$ws.Range("G11").Value2 = 350.3919066666667
$ws.Range("H11").Value2 = 1051.17572
$ws.Range("I11").Value2 = 0.4789344206933965
$ws.Range("J11").Value2 = 0.4937404802104949
$ws.Range("K11").Value2 = 3
$ws.Range("M11").Value2 = 79.34548433333335
$ws.Range("N11").Value2 = 238.036453
$ws.Range("O11").Value2 = 0.2936504262229702
$ws.Range("P11").Value2 = 0.3091729161606711
$ws.Range("Q11").Value2 = 27802.0155409468
$ws.Range("R11").Value2 = 250218.1398685212
$ws.Range("S11").Value2 = 0.1406392967694672
$ws.Range("T11").Value2 = 0.1526511840932488
$ws.Range("E12").Value2 = 3
$ws.Range("G12").Value2 = 350.3919066666667
$ws.Range("H12").Value2 = 1051.17572
$ws.Range("I12").Value2 = 0.4789344206933965
$ws.Range("J12").Value2 = 0.4937404802104949
$ws.Range("K12").Value2 = 3
$ws.Range("M12").Value2 = 11.726012
$ws.Range("N12").Value2 = 35.17803600000001
$ws.Range("O12").Value2 = 0.04339690469630293
$ws.Range("P12").Value2 = 0.04569088405516222
$ws.Range("Q12").Value2 = 4108.699702276214
$ws.Range("R12").Value2 = 36978.29732048592
$ws.Range("S12").Value2 = 0.02078427141061038
$ws.Range("T12").Value2 = 0.02255943903463784
$ws.Range("E13").Value2 = 3
$ws.Range("G13").Value2 = 350.3919066666667
$ws.Range("H13").Value2 = 1051.17572
$ws.Range("I13").Value2 = 0.4789344206933965
$ws.Range("J13").Value2 = 0.4937404802104949
$ws.Range("K13").Value2 = 2
$ws.Range("M13").Value2 = 25.708558
$ws.Range("N13").Value2 = 51.417116
$ws.Range("O13").Value2 = 0.09514503664207198
$ws.Range("P13").Value2 = 0.0667829632560165
$ws.Range("Q13").Value2 = 9008.070655270587
$ws.Range("R13").Value2 = 54048.42393162352
$ws.Range("S13").Value2 = 0.04556823300602272
$ws.Range("T13").Value2 = 0.03297345234790543
$ws.Range("E14").Value2 = 3
$ws.Range("G14").Value2 = 243.8287033333334
$ws.Range("H14").Value2 = 731.4861100000001
$ws.Range("I14").Value2 = 0.3332781281688242
$ws.Range("J14").Value2 = 0.3435812836494235
$ws.Range("K14").Value2 = 2
$ws.Range("M14").Value2 = 14.989415
$ws.Range("N14").Value2 = 29.97883
$ws.Range("O14").Value2 = 0.05547446260572933
$ws.Range("P14").Value2 = 0.03893791130463959
$ws.Range("Q14").Value2 = 3654.849623175217
$ws.Range("R14").Value2 = 21929.0977390513
$ws.Range("S14").Value2 = 0.01848842505840891
$ws.Range("T14").Value2 = 0.01337833754867547
$ws.Range("E15").Value2 = 3
$ws.Range("G15").Value2 = 243.8287033333334
$ws.Range("H15").Value2 = 731.4861100000001
$ws.Range("I15").Value2 = 0.3332781281688242
$ws.Range("J15").Value2 = 0.3435812836494235
$ws.Range("K15").Value2 = 3
$ws.Range("M15").Value2 = 58.255493
$ws.Range("N15").Value2 = 174.766479
$ws.Range("O15").Value2 = 0.2155982850569436
$ws.Range("P15").Value2 = 0.2269949046819425
$ws.Range("Q15").Value2 = 14204.36132023408
$ws.Range("R15").Value2 = 127839.2518821067
$ws.Range("S15").Value2 = 0.07185419288018675
$ws.Range("T15").Value2 = 0.07799120073250032
$ws.Range("E16").Value2 = 3
$ws.Range("G16").Value2 = 243.8287033333334
$ws.Range("H16").Value2 = 731.4861100000001
$ws.Range("I16").Value2 = 0.3332781281688242
$ws.Range("J16").Value2 = 0.3435812836494235
$ws.Range("K16").Value2 = 3
$ws.Range("M16").Value2 = 80.178917
$ws.Range("N16").Value2 = 240.536751
$ws.Range("O16").Value2 = 0.2967348847759819
$ws.Range("P16").Value2 = 0.3124204205415681
$ws.Range("Q16").Value2 = 19549.92136678096
$ws.Range("R16").Value2 = 175949.2923010286
$ws.Range("S16").Value2 = 0.09889524696053099
$ws.Range("T16").Value2 = 0.1073418091279647
$ws.Range("E17").Value2 = 3
$ws.Range("G17").Value2 = 243.8287033333334
$ws.Range("H17").Value2 = 731.4861100000001
$ws.Range("I17").Value2 = 0.3332781281688242
$ws.Range("J17").Value2 = 0.3435812836494235
$ws.Range("K17").Value2 = 3
$ws.Range("M17").Value2 = 79.34548433333335
$ws.Range("N17").Value2 = 238.036453
$ws.Range("O17").Value2 = 0.2936504262229702
$ws.Range("P17").Value2 = 0.3091729161606711
$ws.Range("Q17").Value2 = 19346.70656035198
$ws.Range("R17").Value2 = 174120.3590431678
$ws.Range("S17").Value2 = 0.09786726438756892
$ws.Range("T17").Value2 = 0.106226027404119
$ws.Range("E18").Value2 = 3
$ws.Range("G18").Value2 = 243.8287033333334
$ws.Range("H18").Value2 = 731.4861100000001
$ws.Range("I18").Value2 = 0.3332781281688242
$ws.Range("J18").Value2 = 0.3435812836494235
$ws.Range("K18").Value2 = 3
$ws.Range("M18").Value2 = 11.726012
$ws.Range("N18").Value2 = 35.17803600000001
$ws.Range("O18").Value2 = 0.04339690469630293
$ws.Range("P18").Value2 = 0.04569088405516222
$ws.Range("Q18").Value2 = 2859.138301231108
$ws.Range("R18").Value2 = 25732.24471107997
$ws.Range("S18").Value2 = 0.0144632391655047
$ws.Range("T18").Value2 = 0.01569853259474961
$ws.Range("E19").Value2 = 3
$ws.Range("G19").Value2 = 243.8287033333334
$ws.Range("H19").Value2 = 731.4861100000001
$ws.Range("I19").Value2 = 0.3332781281688242
$ws.Range("J19").Value2 = 0.3435812836494235
$ws.Range("K19").Value2 = 2
$ws.Range("M19").Value2 = 25.708558
$ws.Range("N19").Value2 = 51.417116
$ws.Range("O19").Value2 = 0.09514503664207198
$ws.Range("P19").Value2 = 0.0667829632560165
$ws.Range("Q19").Value2 = 6268.484361709794
$ws.Range("R19").Value2 = 37610.90617025876
$ws.Range("S19").Value2 = 0.03170975971662394
$ws.Range("T19").Value2 = 0.02294537624141443
$ws.Range("E20").Value2 = 3
$ws.Range("G20").Value2 = 48.737294
$ws.Range("H20").Value2 = 146.211882
$ws.Range("I20").Value2 = 0.06661674320651284
$ws.Range("J20").Value2 = 0.06867617226847689
$ws.Range("K20").Value2 = 2
$ws.Range("M20").Value2 = 14.989415
$ws.Range("N20").Value2 = 29.97883
$ws.Range("O20").Value2 = 0.05547446260572933
$ws.Range("P20").Value2 = 0.03893791130463959
$ws.Range("Q20").Value2 = 730.54352574301
$ws.Range("R20").Value2 = 4383.26115445806
$ws.Range("S20").Value2 = 0.00369552802992517
$ws.Range("T20").Value2 = 0.002674106704532102
$ws.Range("E21").Value2 = 3
$ws.Range("G21").Value2 = 48.737294
$ws.Range("H21").Value2 = 146.211882
$ws.Range("I21").Value2 = 0.06661674320651284
$ws.Range("J21").Value2 = 0.06867617226847689
$ws.Range("K21").Value2 = 3
$ws.Range("M21").Value2 = 58.255493
$ws.Range("N21").Value2 = 174.766479
$ws.Range("O21").Value2 = 0.2155982850569436
$ws.Range("P21").Value2 = 0.2269949046819425
$ws.Range("Q21").Value2 = 2839.215089455942
$ws.Range("R21").Value2 = 25552.93580510348
$ws.Range("S21").Value2 = 0.01436245559140297
$ws.Range("T21").Value2 = 0.01558914117800357
$ws.Range("E22").Value2 = 3
$ws.Range("G22").Value2 = 48.737294
$ws.Range("H22").Value2 = 146.211882
$ws.Range("I22").Value2 = 0.06661674320651284
$ws.Range("J22").Value2 = 0.06867617226847689
$ws.Range("K22").Value2 = 3
$ws.Range("M22").Value2 = 80.178917
$ws.Range("N22").Value2 = 240.536751
$ws.Range("O22").Value2 = 0.2967348847759819
$ws.Range("P22").Value2 = 0.3124204205415681
$ws.Range("Q22").Value2 = 3907.703450430598
$ws.Range("R22").Value2 = 35169.33105387539
$ws.Range("S22").Value2 = 0.01976751161953576
$ws.Range("T22").Value2 = 0.02145583862130273
$ws.Range("E23").Value2 = 3
$ws.Range("G23").Value2 = 48.737294
$ws.Range("H23").Value2 = 146.211882
$ws.Range("I23").Value2 = 0.06661674320651284
$ws.Range("J23").Value2 = 0.06867617226847689
$ws.Range("K23").Value2 = 3
$ws.Range("M23").Value2 = 79.34548433333335
$ws.Range("N23").Value2 = 238.036453
$ws.Range("O23").Value2 = 0.2936504262229702
$ws.Range("P23").Value2 = 0.3091729161606711
$ws.Range("Q23").Value2 = 3867.084197526061
$ws.Range("R23").Value2 = 34803.75777773455
$ws.Range("S23").Value2 = 0.01956203503617865
$ws.Range("T23").Value2 = 0.02123281245099761
$ws.Range("E24").Value2 = 3
$ws.Range("G24").Value2 = 48.737294
$ws.Range("H24").Value2 = 146.211882
$ws.Range("I24").Value2 = 0.06661674320651284
$ws.Range("J24").Value2 = 0.06867617226847689
$ws.Range("K24").Value2 = 3
$ws.Range("M24").Value2 = 11.726012
$ws.Range("N24").Value2 = 35.17803600000001
$ws.Range("O24").Value2 = 0.04339690469630293
$ws.Range("P24").Value2 = 0.04569088405516222
$ws.Range("Q24").Value2 = 571.4940942915281
$ws.Range("R24").Value2 = 5143.446848623753
$ws.Range("S24").Value2 = 0.002890960456111123
$ws.Range("T24").Value2 = 0.003137875024471324
$ws.Range("E25").Value2 = 3
$ws.Range("G25").Value2 = 48.737294
$ws.Range("H25").Value2 = 146.211882
$ws.Range("I25").Value2 = 0.06661674320651284
$ws.Range("J25").Value2 = 0.06867617226847689
$ws.Range("K25").Value2 = 2
$ws.Range("M25").Value2 = 25.708558
$ws.Range("N25").Value2 = 51.417116
$ws.Range("O25").Value2 = 0.09514503664207198
$ws.Range("P25").Value2 = 0.0667829632560165
$ws.Range("Q25").Value2 = 1252.965549562052
$ws.Range("R25").Value2 = 7517.793297372312
$ws.Range("S25").Value2 = 0.006338252473359164
$ws.Range("T25").Value2 = 0.004586398289169551
$ws.Range("E26").Value2 = 3
$ws.Range("G26").Value2 = 22.832077
$ws.Range("H26").Value2 = 68.49623099999999
$ws.Range("I26").Value2 = 0.03120810544755168
$ws.Range("J26").Value2 = 0.03217289111905
$ws.Range("K26").Value2 = 2
$ws.Range("M26").Value2 = 14.989415
$ws.Range("N26").Value2 = 29.97883
$ws.Range("O26").Value2 = 0.05547446260572933
$ws.Range("P26").Value2 = 0.03893791130463959
$ws.Range("Q26").Value2 = 342.239477464955
$ws.Range("R26").Value2 = 2053.43686478973
$ws.Range("S26").Value2 = 0.001731252878645864
$ws.Range("T26").Value2 = 0.001252745180807396
$ws.Range("E27").Value2 = 3
$ws.Range("G27").Value2 = 22.832077
$ws.Range("H27").Value2 = 68.49623099999999
$ws.Range("I27").Value2 = 0.03120810544755168
$ws.Range("J27").Value2 = 0.03217289111905
$ws.Range("K27").Value2 = 3
$ws.Range("M27").Value2 = 58.255493
$ws.Range("N27").Value2 = 174.766479
$ws.Range("O27").Value2 = 0.2155982850569436
$ws.Range("P27").Value2 = 0.2269949046819425
$ws.Range("Q27").Value2 = 1330.093901848961
$ws.Range("R27").Value2 = 11970.84511664065
$ws.Range("S27").Value2 = 0.006728414014368403
$ws.Range("T27").Value2 = 0.007303082352911269
$ws.Range("E28").Value2 = 3
$ws.Range("G28").Value2 = 22.832077
$ws.Range("H28").Value2 = 68.49623099999999
$ws.Range("I28").Value2 = 0.03120810544755168
$ws.Range("J28").Value2 = 0.03217289111905
$ws.Range("K28").Value2 = 3
$ws.Range("M28").Value2 = 80.178917
$ws.Range("N28").Value2 = 240.536751
$ws.Range("O28").Value2 = 0.2967348847759819
$ws.Range("P28").Value2 = 0.3124204205415681
$ws.Range("Q28").Value2 = 1830.651206720609
$ws.Range("R28").Value2 = 16475.86086048548
$ws.Range("S28").Value2 = 0.009260533574055943
$ws.Range("T28").Value2 = 0.01005146817345168
$ws.Range("E29").Value2 = 3
$ws.Range("G29").Value2 = 22.832077
$ws.Range("H29").Value2 = 68.49623099999999
$ws.Range("I29").Value2 = 0.03120810544755168
$ws.Range("J29").Value2 = 0.03217289111905
$ws.Range("K29").Value2 = 3
$ws.Range("M29").Value2 = 79.34548433333335
$ws.Range("N29").Value2 = 238.036453
$ws.Range("O29").Value2 = 0.2936504262229702
$ws.Range("P29").Value2 = 0.3091729161606711
$ws.Range("Q29").Value2 = 1811.622207900961
$ws.Range("R29").Value2 = 16304.59987110864
$ws.Range("S29").Value2 = 0.009164273466284949
$ws.Range("T29").Value2 = 0.009946986568596446
$ws.Range("E30").Value2 = 3
$ws.Range("G30").Value2 = 22.832077
$ws.Range("H30").Value2 = 68.49623099999999
$ws.Range("I30").Value2 = 0.03120810544755168
$ws.Range("J30").Value2 = 0.03217289111905
$ws.Range("K30").Value2 = 3
$ws.Range("M30").Value2 = 11.726012
$ws.Range("N30").Value2 = 35.17803600000001
$ws.Range("O30").Value2 = 0.04339690469630293
$ws.Range("P30").Value2 = 0.04569088405516222
$ws.Range("Q30").Value2 = 267.729208886924
$ws.Range("R30").Value2 = 2409.562879982316
$ws.Range("S30").Value2 = 0.001354335177859572
$ws.Range("T30").Value2 = 0.001470007837839872
$ws.Range("E31").Value2 = 3
$ws.Range("G31").Value2 = 22.832077
$ws.Range("H31").Value2 = 68.49623099999999
$ws.Range("I31").Value2 = 0.03120810544755168
$ws.Range("J31").Value2 = 0.03217289111905
$ws.Range("K31").Value2 = 2
$ws.Range("M31").Value2 = 25.708558
$ws.Range("N31").Value2 = 51.417116
$ws.Range("O31").Value2 = 0.09514503664207198
$ws.Range("P31").Value2 = 0.0667829632560165
$ws.Range("Q31").Value2 = 586.9797758149659
$ws.Range("R31").Value2 = 3521.878654889796
$ws.Range("S31").Value2 = 0.002969296336336951
$ws.Range("T31").Value2 = 0.002148601005443336
$ws.Range("E32").Value2 = 2
$ws.Range("G32").Value2 = 25.082339
$ws.Range("H32").Value2 = 50.164678
$ws.Range("I32").Value2 = 0.03428388404538221
$ws.Range("J32").Value2 = 0.02356250409334498
$ws.Range("K32").Value2 = 2
$ws.Range("M32").Value2 = 14.989415
$ws.Range("N32").Value2 = 29.97883
$ws.Range("O32").Value2 = 0.05547446260572933
$ws.Range("P32").Value2 = 0.03893791130463959
$ws.Range("Q32").Value2 = 375.9695884416851
$ws.Range("R32").Value2 = 1503.87835376674
$ws.Range("S32").Value2 = 0.001901880043454716
$ws.Range("T32").Value2 = 0.0009174746945018739
$ws.Range("E33").Value2 = 2
$ws.Range("G33").Value2 = 25.082339
$ws.Range("H33").Value2 = 50.164678
$ws.Range("I33").Value2 = 0.03428388404538221
$ws.Range("J33").Value2 = 0.02356250409334498
$ws.Range("K33").Value2 = 3
$ws.Range("M33").Value2 = 58.255493
$ws.Range("N33").Value2 = 174.766479
$ws.Range("O33").Value2 = 0.2155982850569436
$ws.Range("P33").Value2 = 0.2269949046819425
$ws.Range("Q33").Value2 = 1461.184024038127
$ws.Range("R33").Value2 = 8767.104144228762
$ws.Range("S33").Value2 = 0.007391546605275515
$ws.Range("T33").Value2 = 0.005348568370736722
$ws.Range("E34").Value2 = 2
$ws.Range("G34").Value2 = 25.082339
$ws.Range("H34").Value2 = 50.164678
$ws.Range("I34").Value2 = 0.03428388404538221
$ws.Range("J34").Value2 = 0.02356250409334498
$ws.Range("K34").Value2 = 3
$ws.Range("M34").Value2 = 80.178917
$ws.Range("N34").Value2 = 240.536751
$ws.Range("O34").Value2 = 0.2967348847759819
$ws.Range("P34").Value2 = 0.3124204205415681
$ws.Range("Q34").Value2 = 2011.074776846863
$ws.Range("R34").Value2 = 12066.44866108118
$ws.Range("S34").Value2 = 0.01017322438187962
$ws.Range("T34").Value2 = 0.007361407437855257
$ws.Range("E35").Value2 = 2
$ws.Range("G35").Value2 = 25.082339
$ws.Range("H35").Value2 = 50.164678
$ws.Range("I35").Value2 = 0.03428388404538221
$ws.Range("J35").Value2 = 0.02356250409334498
$ws.Range("K35").Value2 = 3
$ws.Range("M35").Value2 = 79.34548433333335
$ws.Range("N35").Value2 = 238.036453
$ws.Range("O35").Value2 = 0.2936504262229702
$ws.Range("P35").Value2 = 0.3091729161606711
$ws.Range("Q35").Value2 = 1990.170336167856
$ws.Range("R35").Value2 = 11941.02201700714
$ws.Range("S35").Value2 = 0.01006747716250537
$ws.Range("T35").Value2 = 0.007284888102587216
$ws.Range("E36").Value2 = 2
$ws.Range("G36").Value2 = 25.082339
$ws.Range("H36").Value2 = 50.164678
$ws.Range("I36").Value2 = 0.03428388404538221
$ws.Range("J36").Value2 = 0.02356250409334498
$ws.Range("K36").Value2 = 3
$ws.Range("M36").Value2 = 11.726012
$ws.Range("N36").Value2 = 35.17803600000001
$ws.Range("O36").Value2 = 0.04339690469630293
$ws.Range("P36").Value2 = 0.04569088405516222
$ws.Range("Q36").Value2 = 294.1158081020681
$ws.Range("R36").Value2 = 1764.694848612408
$ws.Range("S36").Value2 = 0.001487814448536552
$ws.Range("T36").Value2 = 0.00107659164257831
$ws.Range("E37").Value2 = 2
$ws.Range("G37").Value2 = 25.082339
$ws.Range("H37").Value2 = 50.164678
$ws.Range("I37").Value2 = 0.03428388404538221
$ws.Range("J37").Value2 = 0.02356250409334498
$ws.Range("K37").Value2 = 2
$ws.Range("M37").Value2 = 25.708558
$ws.Range("N37").Value2 = 51.417116
$ws.Range("O37").Value2 = 0.09514503664207198
$ws.Range("P37").Value2 = 0.0667829632560165
$ws.Range("Q37").Value2 = 644.830766957162
$ws.Range("R37").Value2 = 2579.323067828648
$ws.Range("S37").Value2 = 0.003261941403730438
$ws.Range("T37").Value2 = 0.001573573845085596